# Apply the data refresh for data_indiana_covid_deaths_by_date_by_age_group.xlsx
# - bumps a handful of existing running-total C values
# - rewrites the tail of the table (rows 1467-1509) and appends new rows
#   (1510-1518) reflecting newly reported dates/age-groups

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Scattered single-cell increments to the `covid_deaths` (column C) running
#    totals that occur earlier in the sheet.
# ---------------------------------------------------------------------------
$cUpdates = @(
    @(777, 8),
    @(788, 11),
    @(1022, 59),
    @(1082, 66),
    @(1191, 29),
    @(1213, 55),
    @(1249, 19),
    @(1261, 38),
    @(1339, 29),
    @(1340, 2),
    @(1344, 13),
    @(1360, 16),
    @(1410, 2),
    @(1429, 2),
    @(1432, 11),
    @(1436, 5),
    @(1440, 4),
    @(1446, 7),
    @(1449, 9),
    @(1453, 8),
    @(1457, 5)
)

foreach ($upd in $cUpdates) {
    $row = $upd[0]
    $val = $upd[1]
    $ws.Cells.Item($row, 3).Value = $val
}

# ---------------------------------------------------------------------------
# 2) Rewrite the tail of the table: rows 1467-1509 already exist and get new
#    values (the newly-reported day's figures shift everything below it);
#    rows 1510-1518 are brand-new rows appended at the end.
# ---------------------------------------------------------------------------

# A reference cell that already carries the correct date number format
# (custom format "YYYY-MM-DD HH:MM:SS", style index 2 in the original file).
$dateFormat = $ws.Cells.Item(1466, 1).NumberFormat()

$tailRows = @(
    @(1467, 44254, "50-59", 1),
    @(1468, 44254, "60-69", 1),
    @(1469, 44254, "70-79", 2),
    @(1470, 44254, "80+",   4),
    @(1471, 44255, "40-49", 1),
    @(1472, 44255, "60-69", 2),
    @(1473, 44255, "70-79", 5),
    @(1474, 44255, "80+",   7),
    @(1475, 44256, "60-69", 5),
    @(1476, 44256, "70-79", 5),
    @(1477, 44256, "80+",   2),
    @(1478, 44257, "50-59", 2),
    @(1479, 44257, "60-69", 4),
    @(1480, 44257, "70-79", 9),
    @(1481, 44257, "80+",   7),
    @(1482, 44258, "50-59", 1),
    @(1483, 44258, "60-69", 3),
    @(1484, 44258, "70-79", 3),
    @(1485, 44258, "80+",   2),
    @(1486, 44259, "40-49", 1),
    @(1487, 44259, "50-59", 2),
    @(1488, 44259, "60-69", 3),
    @(1489, 44259, "70-79", 3),
    @(1490, 44259, "80+",   2),
    @(1491, 44260, "50-59", 1),
    @(1492, 44260, "60-69", 4),
    @(1493, 44260, "70-79", 3),
    @(1494, 44260, "80+",   5),
    @(1495, 44261, "50-59", 2),
    @(1496, 44261, "60-69", 3),
    @(1497, 44261, "70-79", 3),
    @(1498, 44262, "60-69", 3),
    @(1499, 44262, "80+",   4),
    @(1500, 44263, "40-49", 1),
    @(1501, 44263, "50-59", 1),
    @(1502, 44263, "60-69", 1),
    @(1503, 44263, "70-79", 3),
    @(1504, 44263, "80+",   5),
    @(1505, 44264, "50-59", 1),
    @(1506, 44264, "60-69", 4),
    @(1507, 44264, "70-79", 3),
    @(1508, 44265, "40-49", 1),
    @(1509, 44265, "50-59", 1),
    @(1510, 44265, "70-79", 4),
    @(1511, 44265, "80+",   4),
    @(1512, 44266, "50-59", 1),
    @(1513, 44266, "60-69", 2),
    @(1514, 44266, "70-79", 1),
    @(1515, 44266, "80+",   1),
    @(1516, 44267, "60-69", 2),
    @(1517, 44267, "80+",   3),
    @(1518, 44268, "50-59", 1)
)

foreach ($r in $tailRows) {
    $rowNum = $r[0]
    $dateSerial = $r[1]
    $ageGroup = $r[2]
    $deaths = $r[3]

    $aCell = $ws.Cells.Item($rowNum, 1)
    $aCell.NumberFormat = $dateFormat
    $aCell.Value = $dateSerial

    $ws.Cells.Item($rowNum, 2).Value = $ageGroup
    $ws.Cells.Item($rowNum, 3).Value = $deaths
}
